# Apply row content shift to column A (rows 6-31) and related B/C/D cell
# updates, matching the target diff for "Financial Statement Output.xlsx".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New text values for column A, rows 6 through 31.
$newA = @{
    6  = " inventories"
    7  = " prepaid expenses and other"
    8  = " deposits and other"
    9  = " accounts payable"
    10 = " accrued expenses"
    11 = " deferred rent"
    12 = " other net long-term cash provided liabilities by operating activities"
    13 = " additions to property and equipment"
    14 = " proceeds net from cash sale-leaseback used in investing transactions activities"
    15 = " exercise of employee stock options"
    16 = " tax benefit from stock-based compensation"
    17 = " cash and cash equivalents at end of the period accompanying notes are an integral part of the consolidated financial statements."
    18 = " net cash provided by operating activities"
    19 = " purchases of solar energy systems net of sales"
    20 = " business combinations net of cash acquired"
    21 = " net cash used in investing activities"
    22 = " proceeds from issuances of convertible and other debt"
    23 = " repayments of convertible and other debt"
    24 = " collateralized lease repayments"
    25 = " principal payments on finance leases"
    26 = " debt issuance costs"
    27 = " distributions paid to noncontrolling interests in subsidiaries gl)"
    28 = " payments for buy-outs of noncontrolling interests in subsidiaries"
    29 = " net cash provided by financing activities"
    30 = " net increase in cash and cash equivalents and restricted cash"
    31 = " net increase in cash and cash equivalents and restricted cash"
}

foreach ($r in 6..31) {
    $ws.Cells.Item($r, 1).Value = $newA[$r]
}

# Rows 22-30 previously carried B/C/D figures that belonged to the old
# (pre-shift) row labels; those numbers no longer apply to the relabeled
# rows, so clear them out.
foreach ($r in 22..30) {
    $ws.Range("B" + $r + ":D" + $r).ClearContents()
}

# Row 31 now represents "net increase in cash and cash equivalents and
# restricted cash", picking up the figures that used to sit on that line
# further down the sheet.
$ws.Cells.Item(31, 2).Value = 13118
$ws.Cells.Item(31, 3).Value = 2506
$ws.Cells.Item(31, 4).Value = 312
